$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.151.10"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.586.52"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.65"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.05"
$ws.Range("E6").Value = "  -3.58%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.597.95"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.52"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.049.90"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.075.11"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.57"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.590.08"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.63"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.45"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.18"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.01"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0722"
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("E31").Value = "  -5.81%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.68"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.60"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.68"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("E40").Value = "  -6.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.53"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "272.60"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.82"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0952"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.37"
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.971.17"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.51"
$ws.Range("E51").Value = "  -3.76%  "
